# Generate Report for Handback
# Records that the handback for "a.md" (en-US source) has been processed:
# the file is now in sync with en-US for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$handbackName     = "TestHandback_201702170423"

# ---------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) / F (de-de) show the per-locale
# status for each source file. Row 2 is a.md.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 (a.md) handback is recorded.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("L2").Value = "2017-02-17 08:23:26"
$wsZhCn.Range("M2").Value = $handbackName
$wsZhCn.Range("R2").ClearContents()

# ---------------------------------------------------------------------
# de-de sheet: row 2 (a.md) handback is recorded.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("L2").Value = "2017-02-17 08:23:50"
$wsDeDe.Range("M2").Value = $handbackName
$wsDeDe.Range("R2").ClearContents()

# ---------------------------------------------------------------------
# Column widths grew to fit the new, longer "Handed back..." status
# text and handback name (cosmetic autofit performed by the report
# generator after the values changed).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.17
$wsZhCn.Columns.Item(13).ColumnWidth = 27.17

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.17
$wsDeDe.Columns.Item(13).ColumnWidth = 27.17
